$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-05-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-06 Tuesday", 2) | Out-Null

# Update each division-problem cell in the table individually
# (targeting by row/column avoids ambiguity from duplicate problem text)
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("55÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷4=", 2) | Out-Null

$cell = $t.Cell(1, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("23÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷5=", 2) | Out-Null

$cell = $t.Cell(1, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("73÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷8=", 2) | Out-Null

$cell = $t.Cell(1, 4)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("10÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷6=", 2) | Out-Null

$cell = $t.Cell(1, 5)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("89÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=", 2) | Out-Null

$cell = $t.Cell(5, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("21÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷4=", 2) | Out-Null

$cell = $t.Cell(5, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("75÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷9=", 2) | Out-Null

$cell = $t.Cell(5, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("77÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷8=", 2) | Out-Null

$cell = $t.Cell(5, 4)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("23÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷6=", 2) | Out-Null

$cell = $t.Cell(5, 5)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("73÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷2=", 2) | Out-Null

$cell = $t.Cell(9, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("47÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷7=", 2) | Out-Null

$cell = $t.Cell(9, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("27÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷8=", 2) | Out-Null

$cell = $t.Cell(9, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("14÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=", 2) | Out-Null

$cell = $t.Cell(9, 4)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("60÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷3=", 2) | Out-Null

$cell = $t.Cell(9, 5)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("24÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷6=", 2) | Out-Null

$cell = $t.Cell(13, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("12÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷8=", 2) | Out-Null

$cell = $t.Cell(13, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("46÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷8=", 2) | Out-Null

$cell = $t.Cell(13, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("94÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷9=", 2) | Out-Null

$cell = $t.Cell(13, 4)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("93÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷4=", 2) | Out-Null

$cell = $t.Cell(13, 5)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("17÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=", 2) | Out-Null

$cell = $t.Cell(17, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("64÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷6=", 2) | Out-Null

$cell = $t.Cell(17, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("64÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷8=", 2) | Out-Null

$cell = $t.Cell(17, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("76÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷7=", 2) | Out-Null

$cell = $t.Cell(17, 4)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("51÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷9=", 2) | Out-Null

$cell = $t.Cell(17, 5)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Find.Execute("14÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷5=", 2) | Out-Null
